# Apply the commit's edit:
#   1. "TITLE OF PROJECT" -> "Hospital Equipment Management System"
#   2. Move the (Word-managed) "_GoBack" bookmark from right after the
#      "Tools & Technologies" run to right after the new title text.
#
# Note: this headless engine mis-resolves a zero-length Range that sits
# exactly on a "end-of-paragraph-text / before the paragraph mark"
# boundary when it is handed straight to Bookmarks.Add (it snaps to some
# unrelated earlier position). The reliable workaround is to temporarily
# insert a one-character marker at the target spot (which makes the
# position an ordinary, unambiguous one), add the bookmark right before
# that marker, then delete the marker - the bookmark stays put.

$d = $word.ActiveDocument

# --- 1. Replace the project title text -------------------------------
$found = $d.Content.Find.Execute(
    "TITLE OF PROJECT", $true, $false, $false, $false, $false,
    $true, 1, $false, "Hospital Equipment Management System", 2)
if (-not $found) {
    throw "Could not find 'TITLE OF PROJECT' to replace."
}

# --- 2. Drop the old "_GoBack" bookmark (after Tools & Technologies) -
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 3. Re-create "_GoBack" right after the new title run ------------
$titleRange = $d.Content
$titleRange.Find.Execute("Hospital Equipment Management System")
$pos = $titleRange.End

$marker = $d.Range($pos, $pos)
$marker.InsertAfter("#")

$bmTarget = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmTarget)

$d.Range($pos, $pos + 1).Delete()

Write-Output ("Title updated: " + $found)
Write-Output ("_GoBack now at: " + $d.Bookmarks("_GoBack").Start + "/" + $d.Bookmarks("_GoBack").End)
